$d = $word.ActiveDocument

# 1. asi -> así  (remove spellStart/spellEnd by spanning outside the proofErr-wrapped run)
$d.Content.Find.Execute("datos asi lo", $true, $false, $false, $false, $false, $true, 1, $false, "datos así lo", 2)

# 2. "permite ," -> "permite," (remove gramStart/gramEnd by spanning outside)
$d.Content.Find.Execute("lo permite , esta", $true, $false, $false, $false, $false, $true, 1, $false, "lo permite, esta", 2)

# 3. tecnichal -> Tecnichal (keep its own proofErr pair, narrow match)
$d.Content.Find.Execute("tecnichal", $true, $false, $false, $false, $false, $true, 1, $false, "Tecnichal", 2)

# 4. task -> Task (keep its own proofErr pair, narrow match)
$d.Content.Find.Execute("task", $true, $false, $false, $false, $false, $true, 1, $false, "Task", 2)

# 5. "en el siguiente Sprint" -> "en próximos Sprint"
$d.Content.Find.Execute("debe ser revisada en el siguiente Sprint", $true, $false, $false, $false, $false, $true, 1, $false, "debe ser revisada en próximos Sprint", 2)

# 6. "aclara en en Sprint" -> "aclara en el Sprint" (remove the duplicate "en" proofErr pair by spanning outside)
$d.Content.Find.Execute("aclara en en Sprint", $true, $false, $false, $false, $false, $true, 1, $false, "aclara en el Sprint", 2)
